# Weekly update: insert 3 new price-report rows for "Alcachofa" (artichoke)
# before the existing row 221, pushing the rest of the table down by three
# rows (old A1:R291 -> new A1:R294).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 221..223 (existing rows 221-291 shift to 224-294).
$ws.Range("A221:A223").EntireRow.Insert()

# ---- New row 221 --------------------------------------------------------
$ws.Cells.Item(221, 1).Value = 9
$ws.Cells.Item(221, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(221, 3).Value = "Metropolitana"
$ws.Cells.Item(221, 4).Value = 44463
$ws.Cells.Item(221, 5).Value = 13
$ws.Cells.Item(221, 6).Value = 100112013
$ws.Cells.Item(221, 7).Value = "Alcachofa"
$ws.Cells.Item(221, 8).Value = "Española"
$ws.Cells.Item(221, 9).Value = "Extra"
$ws.Cells.Item(221, 10).Value = 21
$ws.Cells.Item(221, 11).Value = 14000
$ws.Cells.Item(221, 12).Value = 15000
$ws.Cells.Item(221, 13).Value = 14524
$ws.Cells.Item(221, 14).Value = "`$/caja 25 unidades"
$ws.Cells.Item(221, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(221, 16).Value = 14524
$ws.Cells.Item(221, 17).Value = 1
$ws.Cells.Item(221, 18).Value = "Hortaliza"

# ---- New row 222 --------------------------------------------------------
$ws.Cells.Item(222, 1).Value = 9
$ws.Cells.Item(222, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(222, 3).Value = "Metropolitana"
$ws.Cells.Item(222, 4).Value = 44463
$ws.Cells.Item(222, 5).Value = 13
$ws.Cells.Item(222, 6).Value = 100112013
$ws.Cells.Item(222, 7).Value = "Alcachofa"
$ws.Cells.Item(222, 8).Value = "Española"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 52
$ws.Cells.Item(222, 11).Value = 12000
$ws.Cells.Item(222, 12).Value = 13000
$ws.Cells.Item(222, 13).Value = 12500
$ws.Cells.Item(222, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(222, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(222, 16).Value = 417
$ws.Cells.Item(222, 17).Value = 30
$ws.Cells.Item(222, 18).Value = "Hortaliza"

# ---- New row 223 --------------------------------------------------------
$ws.Cells.Item(223, 1).Value = 9
$ws.Cells.Item(223, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(223, 3).Value = "Metropolitana"
$ws.Cells.Item(223, 4).Value = 44463
$ws.Cells.Item(223, 5).Value = 13
$ws.Cells.Item(223, 6).Value = 100112013
$ws.Cells.Item(223, 7).Value = "Alcachofa"
$ws.Cells.Item(223, 8).Value = "Española"
$ws.Cells.Item(223, 9).Value = "Segunda"
$ws.Cells.Item(223, 10).Value = 38
$ws.Cells.Item(223, 11).Value = 10000
$ws.Cells.Item(223, 12).Value = 11000
$ws.Cells.Item(223, 13).Value = 10500
$ws.Cells.Item(223, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(223, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(223, 16).Value = 262
$ws.Cells.Item(223, 17).Value = 40
$ws.Cells.Item(223, 18).Value = "Hortaliza"
